# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" worksheet (fund holdings detail) right after the
# "总计" (Total) summary sheet, pushing "2022-Q3", "2022-Q1" and "2021-Q4"
# down by one position, and adds a corresponding new row to the "总计"
# summary sheet.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q1Sheet = $wb.Worksheets.Item("2022-Q1")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating the "2022-Q1" sheet
#    (same column layout/formatting), placed right after "总计".
# ---------------------------------------------------------------------
$q1Sheet.Copy($null, $totalSheet)
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# Duplicate row 2's formatting into row 3 so the new second data row gets
# the same "A column" bold/bordered style as row 2.
$q4Sheet.Range("A2:H2").Copy()
$q4Sheet.Range("A3").PasteSpecial(-4122)

# --- Row 2: fund 008381 ---
$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2:G2").NumberFormat = "@"
$q4Sheet.Range("B2").Value = "008381"
$q4Sheet.Range("C2").Value = "前海开源新兴产业混合A"
$q4Sheet.Range("D2").Value = "6.29"
$q4Sheet.Range("E2").Value = "90.12"
$q4Sheet.Range("F2").Value = "6.72"
$q4Sheet.Range("G2").Value = "0.4227"
$q4Sheet.Range("B2:G2").ClearFormats()
$q4Sheet.Range("H2").Value = 4

# --- Row 3: fund 014729 ---
$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("B3:G3").NumberFormat = "@"
$q4Sheet.Range("B3").Value = "014729"
$q4Sheet.Range("C3").Value = "前海开源新兴产业混合C"
$q4Sheet.Range("D3").Value = "2.06"
$q4Sheet.Range("E3").Value = "90.12"
$q4Sheet.Range("F3").Value = "6.72"
$q4Sheet.Range("G3").Value = "0.1384"
$q4Sheet.Range("B3:G3").ClearFormats()
$q4Sheet.Range("H3").Value = 4

# ---------------------------------------------------------------------
# 2. Update the "总计" (Total) summary sheet: insert a new row 2 for
#    2022-Q4 and shift the existing quarters down.
# ---------------------------------------------------------------------
$totalSheet.Rows("2").Insert()

# Give the new A2 the same style as the (now shifted) A3 cell.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.5600000000000001

# The shifted-down rows keep their old running-index values in column A;
# renumber them (1, 2, 3) to account for the newly inserted row.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

Write-Output "2022-Q4 sheet added"
